# Add a new "Description" bullet item to the Training table fields list,
# right after the existing "Duration" bullet item (same list: numId=3,
# ListParagraph style, Gill Sans MT 12pt runs).

$d = $word.ActiveDocument

# Locate the "Duration" list item paragraph precisely (exact text match,
# ignoring the trailing paragraph mark) so we don't accidentally match the
# later "The format for Duration is: m." sentence.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "Duration") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Duration' list paragraph"
}

$durationPara = $d.Paragraphs.Item($targetIndex)

# Insert a brand-new paragraph right after "Duration". Word carries the
# paragraph (pPr: ListParagraph style + numPr ilvl/numId) and run (rPr:
# Gill Sans MT, sz 24) formatting over from the paragraph it was split from,
# which matches the existing list items exactly.
$durationPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "Description"
